$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 766302
$ws.Range("E2").Value = 1429167271
$ws.Range("C13").Value = 187831
$ws.Range("E13").Value = 1165047270
$ws.Range("C54").Value = 75190
$ws.Range("E54").Value = 361038094
$ws.Range("C67").Value = 27099
$ws.Range("E67").Value = 168703518
$ws.Range("C69").Value = 17885
$ws.Range("E69").Value = 103673255
$ws.Range("C79").Value = 681
$ws.Range("E79").Value = 20544321
$ws.Range("C81").Value = 88349
$ws.Range("E81").Value = 499611238
$ws.Range("C88").Value = 71256
$ws.Range("E88").Value = 110281502
$ws.Range("C91").Value = 18840
$ws.Range("E91").Value = 75078847
$ws.Range("C96").Value = 29541
$ws.Range("E96").Value = 56427433
$ws.Range("C100").Value = 9328
$ws.Range("E100").Value = 23689436
$ws.Range("C104").Value = 319349
$ws.Range("E104").Value = 561270650
$ws.Range("C112").Value = 145224
$ws.Range("E112").Value = 716129465
$ws.Range("C115").Value = 81798
$ws.Range("E115").Value = 436421192
$ws.Range("C121").Value = 1306091
$ws.Range("E121").Value = 2274443698
$ws.Range("C129").Value = 633285
$ws.Range("E129").Value = 3425593399
$ws.Range("C132").Value = 585553
$ws.Range("E132").Value = 3457723544
$ws.Range("C139").Value = 76633
$ws.Range("E139").Value = 114126711
$ws.Range("C144").Value = 25055
$ws.Range("E144").Value = 92126085
$ws.Range("C145").Value = 70
$ws.Range("E145").Value = 6360066
$ws.Range("C146").Value = 7437
$ws.Range("E146").Value = 37677476
$ws.Range("C151").Value = 39917
$ws.Range("E151").Value = 60356678
$ws.Range("C154").Value = 18426
$ws.Range("E154").Value = 72350486
$ws.Range("C156").Value = 12390
$ws.Range("E156").Value = 39995711
$ws.Range("C159").Value = 43846
$ws.Range("E159").Value = 101311481
$ws.Range("C186").Value = 236812
$ws.Range("E186").Value = 1189632759
$ws.Range("C189").Value = 100464
$ws.Range("E189").Value = 556054990
$ws.Range("C204").Value = 265648
$ws.Range("E204").Value = 1271460352
$ws.Range("C213").Value = 23446
$ws.Range("E213").Value = 88029761
$ws.Range("C215").Value = 230251
$ws.Range("E215").Value = 408697207
$ws.Range("C237").Value = 283304
$ws.Range("E237").Value = 1438285118
$ws.Range("C240").Value = 205892
$ws.Range("E240").Value = 1068191228
